$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(7, 2768.07168,          2709.215670857143,  58.85600914285669,  0.01670676034570443),
    @(2, 3733.67808,          3911.855451428573,  178.1773714285728,  0.005581061894295284),
    @(1, 2446.20288,          2671.051227428571,  224.8483474285713,  0.004427749909997763),
    @(4, 2639.32416,          2982.4363008,       343.1121408000004,  0.002906029405632639),
    @(6, 5391.3024,           5954.342893714289,  563.0404937142894,  0.001772922354235338),
    @(5, 4667.097600000001,   5867.208411428572,  1200.110811428571,  0.0008325626499112312),
    @(0, 4506.1632,           6291.385508571428,  1785.222308571428,  0.0005598407293433554),
    @(3, 14162.2272,          12143.65001142857,  2018.577188571431,  0.0004951531467372932)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
